$words = @(
    "streichen",
    "rufen",
    "wüten",
    "mauern",
    "fahren",
    "knarren",
    "hören",
    "formen",
    "reizen",
    "fischen",
    "bauen",
    "boxen",
    "ändern",
    "sinken",
    "brauchen",
    "kehren",
    "warnen",
    "stehlen",
    "graben",
    "heulen",
    "gelten",
    "spüren",
    "rasen",
    "betteln",
    "helfen",
    "jagen",
    "wenden",
    "klettern",
    "klingen",
    "werden",
    "schwören",
    "weichen",
    "jubeln",
    "scheinen",
    "schreiten",
    "folgen",
    "achten",
    "schlucken",
    "liefern",
    "feiern",
    "quälen",
    "schenken",
    "machen",
    "heben",
    "greifen",
    "wachsen",
    "kümmern",
    "wehtun",
    "bergen",
    "ehren",
    "wirken",
    "spinnen",
    "fallen",
    "schulden",
    "sichern",
    "platzen",
    "sprengen",
    "malen",
    "runden",
    "sterben",
    "schwingen",
    "freuen",
    "grüßen",
    "altern",
    "bellen",
    "saufen",
    "zielen",
    "ärgern",
    "tollen",
    "heilen",
    "äußern",
    "irren",
    "planen",
    "bitten",
    "zünden",
    "lesen",
    "decken",
    "schmecken",
    "werfen",
    "trauen",
    "treiben",
    "töten",
    "dienen",
    "lockern",
    "bluten",
    "loben",
    "zeigen",
    "sorgen",
    "zögern",
    "arten",
    "flüchten",
    "kichern",
    "gründen",
    "pfeifen",
    "suchen",
    "fließen",
    "drehen",
    "scheitern",
    "pflanzen",
    "stecken",
    "kosten",
    "schwächen",
    "wundern",
    "trennen",
    "stammen",
    "backen",
    "erben",
    "münzen",
    "räumen",
    "tropfen",
    "enden",
    "liegen",
    "hauen",
    "führen",
    "siegen",
    "lügen",
    "dringen",
    "seufzen",
    "flehen",
    "sperren",
    "fällen",
    "biegen",
    "filmen",
    "schrecken",
    "fangen",
    "mögen",
    "geben",
    "spielen"
)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($i = 0; $i -lt $words.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 3).Value = $words[$i]
}

